$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 348.14285
$ws.Range("I4").Value = 184.9
$ws.Range("J4").Value = 756.25
$ws.Range("K4").Value = 184.9
$ws.Range("L4").Value = 756.25
$ws.Range("M4").Value = -70.90000000000001
$ws.Range("N4").Value = -984.25

$ws.Range("H33").Value = 7594.143
$ws.Range("I33").Value = 9301.727999999999
$ws.Range("J33").Value = 1333
$ws.Range("K33").Value = 9301.727999999999
$ws.Range("L33").Value = 1333
$ws.Range("M33").Value = -9072.727999999999
$ws.Range("N33").Value = -1791

$ws.Range("H40").Value = 2922.75
$ws.Range("I40").Value = 2922.75
$ws.Range("K40").Value = 2922.75
$ws.Range("M40").Value = -2747.75

$ws.Range("H115").Value = 2224.889
$ws.Range("I115").Value = 3684.6667
$ws.Range("J115").Value = 1495
$ws.Range("K115").Value = 11054.0001
$ws.Range("L115").Value = 4485
$ws.Range("M115").Value = -9487.000100000001
$ws.Range("N115").Value = -7619

$ws.Range("H137").Value = 1899.6428
$ws.Range("I137").Value = 1507.9166
$ws.Range("K137").Value = 4523.7498
$ws.Range("M137").Value = -1973.7498

$ws.Range("H138").Value = 2506.8772
$ws.Range("I138").Value = 2248.7856
$ws.Range("J138").Value = 3229.5334
$ws.Range("K138").Value = 6746.3568
$ws.Range("L138").Value = 9688.600199999999
$ws.Range("M138").Value = -1606.3568
$ws.Range("N138").Value = -19968.6002


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 14287.375
$ws.Range("I28").Value = 14287.375
$ws.Range("K28").Value = 14287.375
$ws.Range("M28").Value = -14095.375

$ws.Range("H31").Value = 7005.143
$ws.Range("I31").Value = 7005.143
$ws.Range("K31").Value = 7005.143
$ws.Range("M31").Value = -6711.143

$ws.Range("H53").Value = 9874.25
$ws.Range("I53").Value = 5570.7144
$ws.Range("K53").Value = 5570.7144
$ws.Range("M53").Value = -4888.7144

$ws.Range("H99").Value = 14287.375
$ws.Range("I99").Value = 14287.375
$ws.Range("K99").Value = 14287.375
$ws.Range("M99").Value = -11292.375

$ws.Range("H104").Value = 64833
$ws.Range("J104").Value = 64833
$ws.Range("L104").Value = 64833
$ws.Range("N104").Value = -71821

$ws.Range("H121").Value = 38750
$ws.Range("J121").Value = 38750
$ws.Range("L121").Value = 38750
$ws.Range("N121").Value = -42244

$ws.Range("H123").Value = 90832.836
$ws.Range("I123").Value = 89000
$ws.Range("J123").Value = 99997
$ws.Range("K123").Value = 89000
$ws.Range("L123").Value = 99997
$ws.Range("M123").Value = -84100
$ws.Range("N123").Value = -109797

$ws.Range("H132").Value = 17706.213
$ws.Range("I132").Value = 21777.654
$ws.Range("K132").Value = 65332.962
$ws.Range("M132").Value = -62802.962

$ws.Range("H139").Value = 79949
$ws.Range("J139").Value = 79949
$ws.Range("L139").Value = 79949
$ws.Range("N139").Value = -90229


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 10000
$ws.Range("L56").Value = 10000
$ws.Range("N56").Value = -11478

$ws.Range("H107").Value = 1910.9231
$ws.Range("I107").Value = 1485.6364
$ws.Range("K107").Value = 1485.6364
$ws.Range("M107").Value = 434.3635999999999


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3394
$ws.Range("I31").Value = 1669.75
$ws.Range("J31").Value = 4985.615
$ws.Range("K31").Value = 1669.75
$ws.Range("L31").Value = 4985.615
$ws.Range("M31").Value = -1374.75
$ws.Range("N31").Value = -5575.615

$ws.Range("H34").Value = 3394
$ws.Range("I34").Value = 1669.75
$ws.Range("J34").Value = 4985.615
$ws.Range("K34").Value = 1669.75
$ws.Range("L34").Value = 4985.615
$ws.Range("M34").Value = -1467.75
$ws.Range("N34").Value = -5389.615

$ws.Range("H58").Value = 52322.6
$ws.Range("I58").Value = 73435.71000000001
$ws.Range("K58").Value = 73435.71000000001
$ws.Range("M58").Value = -73232.71000000001

$ws.Range("H69").Value = 79999
$ws.Range("J69").Value = 79999
$ws.Range("L69").Value = 79999
$ws.Range("N69").Value = -81497

$ws.Range("H72").Value = 79999
$ws.Range("J72").Value = 79999
$ws.Range("L72").Value = 239997
$ws.Range("N72").Value = -247485

$ws.Range("H99").Value = 3610.1177
$ws.Range("I99").Value = 3539.4167
$ws.Range("J99").Value = 3779.8
$ws.Range("K99").Value = 3539.4167
$ws.Range("L99").Value = 3779.8
$ws.Range("M99").Value = -2041.4167
$ws.Range("N99").Value = -6775.8

$ws.Range("H105").Value = 803
$ws.Range("I105").Value = 803
$ws.Range("K105").Value = 803
$ws.Range("M105").Value = 944

$ws.Range("H126").Value = 3610.1177
$ws.Range("I126").Value = 3539.4167
$ws.Range("J126").Value = 3779.8
$ws.Range("K126").Value = 10618.2501
$ws.Range("L126").Value = 11339.4
$ws.Range("M126").Value = -8148.250100000001
$ws.Range("N126").Value = -16279.4

$ws.Range("H136").Value = 52322.6
$ws.Range("I136").Value = 73435.71000000001
$ws.Range("K136").Value = 220307.13
$ws.Range("M136").Value = -217757.13


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null

$ws.Range("H36").Value = 19999
$ws.Range("J36").Value = 19999
$ws.Range("L36").Value = 19999
$ws.Range("N36").Value = -20969

$ws.Range("H97").Value = 1429.5714
$ws.Range("I97").Value = 1687.6666
$ws.Range("J97").Value = 1236
$ws.Range("K97").Value = 1687.6666
$ws.Range("L97").Value = 1236
$ws.Range("M97").Value = -1191.6666
$ws.Range("N97").Value = -2228

$ws.Range("H102").Value = 4117.222
$ws.Range("I102").Value = 3722.2856
$ws.Range("K102").Value = 3722.2856
$ws.Range("M102").Value = -2100.2856

$ws.Range("H126").Value = 5537.4
$ws.Range("I126").Value = 4259.875
$ws.Range("J126").Value = 6997.4287
$ws.Range("K126").Value = 12779.625
$ws.Range("L126").Value = 20992.2861
$ws.Range("M126").Value = -10309.625
$ws.Range("N126").Value = -25932.2861


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5799.4
$ws.Range("I40").Value = 5305.5386
$ws.Range("K40").Value = 5305.5386
$ws.Range("M40").Value = -5169.5386

$ws.Range("H46").Value = 2395.724
$ws.Range("I46").Value = 1750.5
$ws.Range("J46").Value = 2443.5186
$ws.Range("K46").Value = 1750.5
$ws.Range("L46").Value = 2443.5186
$ws.Range("M46").Value = -1562.5
$ws.Range("N46").Value = -2819.5186

$ws.Range("H61").Value = 4145.8184
$ws.Range("I61").Value = 3710.45
$ws.Range("K61").Value = 3710.45
$ws.Range("M61").Value = -3508.45

$ws.Range("H93").Value = 1467.625
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null

$ws.Range("H113").Value = 4145.8184
$ws.Range("I113").Value = 3710.45
$ws.Range("K113").Value = 3710.45
$ws.Range("M113").Value = -1540.45

$ws.Range("H122").Value = 4622.5
$ws.Range("I122").Value = 4030.7144
$ws.Range("K122").Value = 12092.1432
$ws.Range("M122").Value = -9642.143199999999

$ws.Range("H132").Value = 64067.5
$ws.Range("J132").Value = 5812.222
$ws.Range("L132").Value = 17436.666
$ws.Range("N132").Value = -22496.666


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = $null

$ws.Range("H107").Value = 2086.7222
$ws.Range("I107").Value = 1186.7
$ws.Range("K107").Value = 3560.1
$ws.Range("M107").Value = -1640.1

$ws.Range("H122").Value = 1979.7778
$ws.Range("J122").Value = 2200
$ws.Range("L122").Value = 6600
$ws.Range("N122").Value = -11500

